$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.593.95"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "'2.527.31"
$ws.Range("E3").Value = "  -5.35%  "
$ws.Range("D5").Value = "'575.42"
$ws.Range("E5").Value = "  -4.03%  "
$ws.Range("D6").Value = "'169.63"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.507"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").Value = "'2.525.15"
$ws.Range("E9").Value = "  -5.41%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").Value = "'4.80"
$ws.Range("E13").Value = "  -3.87%  "
$ws.Range("D14").Value = "'2.985.11"
$ws.Range("E14").Value = "  -5.89%  "
$ws.Range("D15").Value = "'70.342.76"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "'24.88"
$ws.Range("E17").Value = "  -5.06%  "
$ws.Range("D18").Value = "'2.529.67"
$ws.Range("E18").Value = "  -5.44%  "
$ws.Range("D19").Value = "'11.52"
$ws.Range("E19").Value = "  -5.74%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("E20").Value = "  -7.75%  "
$ws.Range("D21").Value = "'356.57"
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("D22").Value = "'3.93"
$ws.Range("E22").Value = "  -5.91%  "
$ws.Range("D23").Value = "'1.96"
$ws.Range("E23").Value = "  -3.99%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'69.20"
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("D26").Value = "'4.06"
$ws.Range("E26").Value = "  -6.46%  "
$ws.Range("D27").Value = "'9.20"
$ws.Range("E27").Value = "  -6.09%  "
$ws.Range("D28").Value = "'2.648.53"
$ws.Range("E28").Value = "  -5.87%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "'0.0₃0911"
$ws.Range("E30").Value = "  -6.32%  "
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("D32").Value = "'480.09"
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("D33").Value = "'1.27"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").Value = "'1.76"
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'157.27"
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("D37").Value = "'0.116"
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("D38").Value = "'18.84"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").Value = "'18.56"
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("E42").Value = "  -7.06%  "
$ws.Range("D43").Value = "'0.319"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("D44").Value = "'4.71"
$ws.Range("E44").Value = "  -5.62%  "
$ws.Range("D45").Value = "'2.41"
$ws.Range("E45").Value = "  -6.03%  "
$ws.Range("D46").Value = "'38.31"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Value = "'142.58"
$ws.Range("E47").Value = "  -8.64%  "
$ws.Range("D48").Value = "'3.53"
$ws.Range("E48").Value = "  -5.48%  "
$ws.Range("D49").Value = "'0.524"
$ws.Range("E49").Value = "  -6.44%  "
$ws.Range("E50").Value = "  -6.95%  "
$ws.Range("D51").Value = "'0.596"
$ws.Range("E51").Value = "  -1.53%  "
